$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = -4
$ws.Range("F6").Value = -9
$ws.Range("F8").Value = -8
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = -7
$ws.Range("F13").Value = -6
$ws.Range("F14").Value = -5
$ws.Range("F15").Value = 0
$ws.Range("F19").Value = -3
$ws.Range("F22").Value = -6
$ws.Range("F23").Value = -4
$ws.Range("F24").Value = -4
